# Updated cryptos list on Tue Aug 22 02:37:00 UTC 2023 with GitHub Actions
#
# The "Price" column (D) holds numbers that are stored as plain text
# (e.g. "26.144.98", "1.670.11") because they use '.' as a thousands
# separator rather than a decimal point. Some of the updated prices
# (e.g. "1.000", "210.85") *do* parse as valid numbers, so Excel would
# silently coerce them to numeric cells on a plain Value assignment.
# Forcing the column to the Text ("@") number format before writing
# keeps every cell in column D a literal string, matching the rest of
# the column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.131.98"
$ws.Range("E2").Value = "  -0.75%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.670.71"
$ws.Range("E3").Value = "  -1.39%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.64%  "

# Row 5 - BNB
$ws.Range("D5").Value = "210.85"
$ws.Range("E5").Value = "  -3.26%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.5264"
$ws.Range("E6").Value = "  -2.73%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.62%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -3.66%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.06315"
$ws.Range("E9").Value = "  -2.19%  "

# Row 10 - Solana
$ws.Range("D10").Value = "21.21"
$ws.Range("E10").Value = "  -2.04%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.07559"
$ws.Range("E11").Value = "  -1.27%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.679.47"
$ws.Range("E12").Value = "  -2.83%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "4.444"
$ws.Range("E13").Value = "  -2.03%  "

# Row 14 - Polygon
$ws.Range("E14").Value = "  -3.91%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "67.06"
$ws.Range("E15").Value = "  +0.11%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "0.000007968"
$ws.Range("E16").Value = "  -5.17%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.171.06"
$ws.Range("E17").Value = "  -0.85%  "

# Row 18 - Dai
$ws.Range("E18").Value = "  -0.60%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "4.751"
$ws.Range("E19").Value = "  -3.37%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "187.06"
$ws.Range("E20").Value = "  -1.96%  "

# Row 21 - Avalanche
$ws.Range("E21").Value = "  -4.59%  "

# Row 22 - Chainlink
$ws.Range("D22").Value = "6.195"
$ws.Range("E22").Value = "  -1.38%  "

# Row 23 - BinanceUSD
$ws.Range("E23").Value = "  -0.58%  "

# Row 24 - Monero
$ws.Range("D24").Value = "149.59"

# Row 25 - Stellar
$ws.Range("E25").Value = "  -2.58%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "7.518"
$ws.Range("E26").Value = "  -4.30%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  +0.85%  "

# Row 28 - Hedera
$ws.Range("D28").Value = "0.06276"
$ws.Range("E28").Value = "  -0.50%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  -2.12%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "1.283"
$ws.Range("E30").Value = "  -3.33%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "3.509"
$ws.Range("E31").Value = "  -2.83%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "3.420"
$ws.Range("E32").Value = "  -4.76%  "

# Row 33 - LidoDAOToken
$ws.Range("D33").Value = "1.631"
$ws.Range("E33").Value = "  -3.35%  "

# Row 34 - ARBITRUM
$ws.Range("D34").Value = "1.000"
$ws.Range("E34").Value = "  -3.12%  "

# Row 35 - ImmutableX
$ws.Range("D35").Value = "0.6056"
$ws.Range("E35").Value = "  -2.26%  "

# Row 36 - HuobiToken
$ws.Range("D36").Value = "2.410"
$ws.Range("E36").Value = "  -0.39%  "

# Row 37 - MXToken
$ws.Range("D37").Value = "2.736"
$ws.Range("E37").Value = "  -1.39%  "

# Row 38 - FraxShare
$ws.Range("D38").Value = "6.128"
$ws.Range("E38").Value = "  +0.43%  "

# Row 39/40 - Maker and VeChain swap places (VeChain now ranks 39, Maker 40)
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01612"
$ws.Range("E39").Value = "  -2.54%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.100.34"
$ws.Range("E40").Value = "  -1.29%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "0.8752"
$ws.Range("E41").Value = "  -1.33%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -1.08%  "

# Row 43 - Quant
$ws.Range("D43").Value = "100.27"
$ws.Range("E43").Value = "  -0.88%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.822.84"
$ws.Range("E44").Value = "  -1.30%  "

# Row 45 - BabyDogeCoin
$ws.Range("D45").Value = "0.00000000113"
$ws.Range("E45").Value = "  +1.21%  "

# Row 46 - Aave
$ws.Range("D46").Value = "55.47"
$ws.Range("E46").Value = "  -3.95%  "

# Row 47 - Frax
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  -0.26%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "8.046"
$ws.Range("E48").Value = "  -1.44%  "

# Row 49 - Cronos
$ws.Range("D49").Value = "0.05237"
$ws.Range("E49").Value = "  -0.94%  "

# Row 50 - Mantle
$ws.Range("D50").Value = "0.4246"
$ws.Range("E50").Value = "  -1.32%  "

# Row 51 - Aptos
$ws.Range("D51").Value = "5.998"
$ws.Range("E51").Value = "  -1.64%  "
